## SubscriptionExpirySample.xlsx edit
## - ICCID value: drop trailing space
## - Expiry Date cell: convert from a date serial to plain text "24/10/2026"
##   and give the column a Text number format
## - Remove the empty trailing column C (dimension becomes A1:B2)
## - Re-style header row (bold, 12pt, "ariel", dark grey) and body row
##   (10pt Arial) to match the new layout
## - Re-size columns / rows, update selection + header/footer/margins

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Data -----------------------------------------------------------
$ws.Range("A2").Value = "8991102105546012952F"

$ws.Range("B1:B2").NumberFormat = "@"
$ws.Range("B2").Value = "24/10/2026"

$ws.Range("A1:A2").NumberFormat = "General"

# ---- Remove the stray empty column C --------------------------------
$ws.Columns.Item(3).Delete()

# ---- Column widths (character units; raw stored width = value + 5/6) -
$ws.Columns.Item(1).ColumnWidth = 22.17
$ws.Columns.Item(2).ColumnWidth = 32.505

# ---- Row heights ------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 12.8

# ---- Fonts --------------------------------------------------------
# Header row: bold, 12pt, "ariel", dark grey, family 0
$ws.Range("A1:B1").Font.Name = "ariel"
$ws.Range("A1:B1").Font.Family = 0
$ws.Range("A1:B1").Font.Size = 12
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Color = 3355443

# Body row: Arial 10pt (plain)
$ws.Range("A2:B2").Font.Name = "Arial"
$ws.Range("A2:B2").Font.Size = 10

# ---- Selection ------------------------------------------------------
[void]$ws.Range("A1:B1").Select()

# ---- Page setup / margins / header-footer ----------------------------
$ws.PageSetup.LeftMargin = 56.7
$ws.PageSetup.RightMargin = 56.7
$ws.PageSetup.TopMargin = 75.8
$ws.PageSetup.BottomMargin = 75.8
$ws.PageSetup.HeaderMargin = 56.7
$ws.PageSetup.FooterMargin = 56.7
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'
